# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (per commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 2132   # was 2128
$ws1.Range("F4").Value  = 19     # was 17
$ws1.Range("F5").Value  = 11093  # was 11086
$ws1.Range("F6").Value  = 192    # was 191
$ws1.Range("F7").Value  = 168    # was 167
$ws1.Range("F10").Value = 10987  # was 10977
$ws1.Range("F15").Value = 5504   # was 5499
$ws1.Range("F17").Value = 3424   # was 3422

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 2132   # was 2128
$ws4.Range("F5").Value  = 19     # was 17
$ws4.Range("F7").Value  = 11093  # was 11086
$ws4.Range("F8").Value  = 192    # was 191
$ws4.Range("F9").Value  = 168    # was 167
$ws4.Range("F12").Value = 10987  # was 10977
$ws4.Range("F17").Value = 5504   # was 5499
$ws4.Range("F19").Value = 3424   # was 3422
